$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 12, pushing the existing weekly
# records (rows 12-33) down to rows 13-34.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with this week's record. It repeats
# the same market/product metadata as the row that used to be at 12
# (now at 13), but carries its own date and price figures.
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = "2022-03-25"
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 100112052
$ws.Range("G12").Value = "Albahaca"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 1200
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1350
$ws.Range("N12").Value = "`$/paquete"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 1350
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"
